$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: add note in column D
$ws.Range("D4").Value = "program is coded to allow new audio but loop is still in progress"

# Row 21: change status from Pending to Resolved, and add note in column D
$ws.Range("C21").Value = "Resolved"
$ws.Range("D21").Value = "cannot select middle low 3rd or middle low 5th"

# Update the selection to D21
$ws.Range("D21").Select()
